# Updates the cryptocurrency price ("Price", column D) and volume change
# ("Volume(1h)", column E) figures for rows 2-51 to reflect the latest
# scrape, as captured by the automated "Updated cryptos list" GitHub
# Actions commit. NumberFormat is forced to Text ("@") immediately before
# each write so that values such as "1.001" or "  -0.51%  " are stored
# as literal text instead of being auto-converted to numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.548.04"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.873.92"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.42"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4750"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2906"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06490"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.95"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07752"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7376"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.874.89"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "95.97"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.176"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "273.90"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.613.27"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.20"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007483"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.120.45"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.206"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.170"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.178"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.83"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.77"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.906"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09879"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.344"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.502"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.261"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.084"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04776"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.118"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6933"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.720"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.759"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.272"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.29"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.972"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4176"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8357"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.51"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.355"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.27"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.957"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "917.05"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05669"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.93%  "
